# Pruebas de integridad completas
# Refresh the run timestamp (and derived run_id) metadata that gets stamped
# on every evaluation run of the report.

$wb = $excel.ActiveWorkbook

$oldTimestamp = "2025-10-26T10:34:41.077661"
$newTimestamp = "2025-10-26T11:24:53.518499"

$oldRunId = "3009e1c516db48e99fd3d706547deac8"
$newRunId = "32a45db9b48e454bb75cb54551a0162d"

# --- Sheet "preguntas": refresh the per-row "timestamp" column (T2:T23) ---
$wsPreguntas = $wb.Worksheets.Item("preguntas")

$lastRow = $wsPreguntas.Cells.Item($wsPreguntas.Rows.Count, 2).End(-4162).Row
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $wsPreguntas.Range("T$row")
    if ($cell.Text -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}

# --- Sheet "indice_global": refresh run_id, timestamp, and extra_config ---
$wsIndice = $wb.Worksheets.Item("indice_global")

$wsIndice.Range("L2").Value = $newRunId
$wsIndice.Range("Q2").Value = $newTimestamp

$extraConfig = $wsIndice.Range("R2").Text
$wsIndice.Range("R2").Value = $extraConfig.Replace($oldRunId, $newRunId)
